$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1) F12 / F13 : replace hard-coded price with a formula (1.44 * Qty)
# ------------------------------------------------------------------
$ws.Range("F12").Formula = "=1.44*E12"
$ws.Range("F13").Formula = "=1.44*E13"

# ------------------------------------------------------------------
# 2) E26 : quantity was blank, now set to 1
# ------------------------------------------------------------------
$ws.Range("E26").Value = 1

# ------------------------------------------------------------------
# 3) Insert a brand new "Miscellaneous" line item in row 27
#    (Temperature sensor / TMP36 / Adafruit/Sparkfun)
#    First copy formatting from the most similar existing row (13)
#    so the new row blends in with the rest of the table, then
#    overwrite the values.
# ------------------------------------------------------------------
$ws.Range("A13:E13").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)
$ws.Range("G13:H13").Copy()
$ws.Range("G27:H27").PasteSpecial(-4122)

$ws.Range("A27").Value = $null
$ws.Range("B27").Value = "Temperature sensor"
$ws.Range("C27").Value = "TMP36"
$ws.Range("D27").Value = "Adafruit/Sparkfun"
$ws.Range("E27").Value = 1
$ws.Range("G27").Value = "X"
$ws.Range("H27").Value = $null

# F27 needs its own currency format (matches builtin format 8) plus the
# same border as the rest of the table.
$ws.Range("F13").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Value = 1.5
$ws.Range("F27").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# ------------------------------------------------------------------
# 4) Totals: "Required Materials" (row 27->29) now also includes F27,
#    and "Required and Optional Materials" (row 29->30) extends to F27
# ------------------------------------------------------------------
$ws.Range("B29").Formula = "=SUM(F3:F22) + F27"
$ws.Range("B30").Formula = "=SUM(F3:F27)"

# ------------------------------------------------------------------
# 5) Update the view: scroll down one row and move the selection from
#    B30 to C30 (matches the saved workbook view of the edited file)
# ------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("C30").Select()
